$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B2").Value = "C01号直流"
$ws.Range("C2").Value = "2025-01-25 13:46:36"
$ws.Range("D2").Value = 45986.334270833337

$ws.Range("A3").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B3").Value = "C02号直流"
$ws.Range("C3").Value = "2025-01-25 17:13:47"
$ws.Range("D3").Value = 45986.334270833337

$ws.Range("A4").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B4").Value = "C03号直流"
$ws.Range("C4").Value = "2025-01-25 14:14:24"
$ws.Range("D4").Value = 45986.334270833337

$ws.Range("A5").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B5").Value = "C04号直流"
$ws.Range("C5").Value = "2025-01-25 06:24:40"
$ws.Range("D5").Value = 45986.334270833337

$ws.Range("A6").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B6").Value = "C05号直流"
$ws.Range("C6").Value = "2025-01-25 16:01:40"
$ws.Range("D6").Value = 45986.334270833337

$ws.Range("A7").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B7").Value = "D01号直流"
$ws.Range("C7").Value = "2025-01-25 18:30:24"
$ws.Range("D7").Value = 45986.334270833337

$ws.Range("A8").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B8").Value = "D02号直流"
$ws.Range("C8").Value = "2025-01-25 15:39:19"
$ws.Range("D8").Value = 45986.334270833337

$ws.Range("A9").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B9").Value = "D03号直流"
$ws.Range("C9").Value = "2025-01-25 16:09:35"
$ws.Range("D9").Value = 45986.334270833337

$ws.Range("A10").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B10").Value = "D04号直流"
$ws.Range("C10").Value = "2025-01-25 18:29:02"
$ws.Range("D10").Value = 45986.334270833337

$ws.Range("A11").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B11").Value = "D05号直流"
$ws.Range("C11").Value = "2025-01-25 18:27:29"
$ws.Range("D11").Value = 45986.334270833337

$ws.Range("A12").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B12").Value = "E01号直流"
$ws.Range("C12").Value = "2025-01-25 15:22:58"
$ws.Range("D12").Value = 45986.334270833337

$ws.Range("A13").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B13").Value = "E02号直流"
$ws.Range("C13").Value = "2025-01-25 16:45:57"
$ws.Range("D13").Value = 45986.334270833337

$ws.Range("A14").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B14").Value = "E03号直流"
$ws.Range("C14").Value = "2025-01-25 02:54:59"
$ws.Range("D14").Value = 45986.334270833337

$ws.Range("A15").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B15").Value = "E04号直流"
$ws.Range("C15").Value = "2025-01-25 17:08:37"
$ws.Range("D15").Value = 45986.334270833337

$ws.Range("A16").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B16").Value = "004B号直流"
$ws.Range("C16").Value = "2025-02-19 00:26:27"
$ws.Range("D16").Value = 45986.334270833337

$ws.Range("A17").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B17").Value = "701号直流"
$ws.Range("C17").Value = 45927.457337962966
$ws.Range("D17").Value = 45986.334270833337

$ws.Range("A18").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B18").Value = "101号直流"
$ws.Range("C18").Value = 45979.18986111111
$ws.Range("D18").Value = 45986.334270833337

$ws.Range("A19").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B19").Value = "603号直流"
$ws.Range("C19").Value = 45980.250173611108
$ws.Range("D19").Value = 45986.334270833337

$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "505号直流"
$ws.Range("C20").Value = 45982.551504629628
$ws.Range("D20").Value = 45986.334270833337

$ws.Range("A21").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value = "502号直流"
$ws.Range("C21").Value = 45982.555462962962
$ws.Range("D21").Value = 45986.334270833337

$ws.Range("A22").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B22").Value = "406号直流"
$ws.Range("C22").Value = 45982.584861111114
$ws.Range("D22").Value = 45986.334270833337

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "702号直流"
$ws.Range("C23").Value = 45983.211712962962
$ws.Range("D23").Value = 45986.334270833337

$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "705号直流"
$ws.Range("C24").Value = 45984.586273148147
$ws.Range("D24").Value = 45986.334270833337

$ws.Range("A25").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B25").Value = "111号直流"
$ws.Range("C25").Value = 45984.66097222222
$ws.Range("D25").Value = 45986.334270833337

$ws.Range("A26").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B26").Value = "B02号直流"
$ws.Range("C26").Value = 45985.02648148148
$ws.Range("D26").Value = 45986.334270833337

$ws.Range("A27").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B27").Value = "202号直流"
$ws.Range("C27").Value = 45985.041481481479
$ws.Range("D27").Value = 45986.334270833337

$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "B03号直流"
$ws.Range("C28").Value = 45985.047025462962
$ws.Range("D28").Value = 45986.334270833337

$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "102号直流"
$ws.Range("C29").Value = 45985.049085648148
$ws.Range("D29").Value = 45986.334270833337

$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "A01号直流"
$ws.Range("C30").Value = 45985.165578703702
$ws.Range("D30").Value = 45986.334270833337

$ws.Range("A31").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B31").Value = "306号直流"
$ws.Range("C31").Value = 45985.175949074073
$ws.Range("D31").Value = 45986.334270833337

$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "602号直流"
$ws.Range("C32").Value = 45985.251747685186
$ws.Range("D32").Value = 45986.334270833337

$ws.Range("A33").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B33").Value = "206号直流"
$ws.Range("C33").Value = 45985.41982638889
$ws.Range("D33").Value = 45986.334270833337

$ws.Range("A34").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B34").Value = "210号直流"
$ws.Range("C34").Value = 45985.461655092593
$ws.Range("D34").Value = 45986.334270833337

$ws.Range("A35").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B35").Value = "108号直流"
$ws.Range("C35").Value = 45985.518576388888
$ws.Range("D35").Value = 45986.334270833337

$ws.Range("A36").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B36").Value = "110号直流"
$ws.Range("C36").Value = 45985.530405092592
$ws.Range("D36").Value = 45986.334270833337

$ws.Range("A37").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B37").Value = "903号直流"
$ws.Range("C37").Value = 45985.546990740739
$ws.Range("D37").Value = 45986.334270833337

$ws.Range("A38").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B38").Value = "905号直流"
$ws.Range("C38").Value = 45985.554803240739
$ws.Range("D38").Value = 45986.334270833337

$ws.Range("A39").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B39").Value = "102号直流"
$ws.Range("C39").Value = 45985.555949074071
$ws.Range("D39").Value = 45986.334270833337

$ws.Range("A40").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B40").Value = "801号直流"
$ws.Range("C40").Value = 45985.558240740742
$ws.Range("D40").Value = 45986.334270833337

$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "B01号直流"
$ws.Range("C41").Value = 45985.559560185182
$ws.Range("D41").Value = 45986.334270833337

$ws.Range("A42").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B42").Value = "204号直流"
$ws.Range("C42").Value = 45985.561053240737
$ws.Range("D42").Value = 45986.334270833337

$ws.Range("A43").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B43").Value = "901号直流"
$ws.Range("C43").Value = 45985.565428240741
$ws.Range("D43").Value = 45986.334270833337

$ws.Range("A44").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B44").Value = "904号直流"
$ws.Range("C44").Value = 45985.569664351853
$ws.Range("D44").Value = 45986.334270833337

$ws.Range("A45").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B45").Value = "604号直流"
$ws.Range("C45").Value = 45985.570324074077
$ws.Range("D45").Value = 45986.334270833337

$ws.Range("A46").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B46").Value = "301号直流"
$ws.Range("C46").Value = 45985.573148148149
$ws.Range("D46").Value = 45986.334270833337

$ws.Range("A47").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B47").Value = "106号直流"
$ws.Range("C47").Value = 45985.578043981484
$ws.Range("D47").Value = 45986.334270833337

$ws.Range("A48").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B48").Value = "109号直流"
$ws.Range("C48").Value = 45985.604872685188
$ws.Range("D48").Value = 45986.334270833337

$ws.Range("A49").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B49").Value = "404号直流"
$ws.Range("C49").Value = 45985.627187500002
$ws.Range("D49").Value = 45986.334270833337

$ws.Range("A50").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B50").Value = "105号直流"
$ws.Range("C50").Value = 45985.638564814813
$ws.Range("D50").Value = 45986.334270833337

$ws.Range("A51").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B51").Value = "004A号直流"
$ws.Range("C51").Value = 45985.666342592594
$ws.Range("D51").Value = 45986.334270833337

$ws.Range("A52").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B52").Value = "404号直流"
$ws.Range("C52").Value = 45985.704155092593
$ws.Range("D52").Value = 45986.334270833337

$ws.Range("A53").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B53").Value = "101号直流"
$ws.Range("C53").Value = 45985.724363425928
$ws.Range("D53").Value = 45986.334270833337

$ws.Range("A54").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B54").Value = "006B号直流"
$ws.Range("C54").Value = 45985.735844907409
$ws.Range("D54").Value = 45986.334270833337

$ws.Range("A55").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B55").Value = "305号直流"
$ws.Range("C55").Value = 45985.741284722222
$ws.Range("D55").Value = 45986.334270833337

$ws.Range("A56").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B56").Value = "A03号直流"
$ws.Range("C56").Value = 45985.776724537034
$ws.Range("D56").Value = 45986.334270833337

$ws.Range("A57").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B57").Value = "502号直流"
$ws.Range("C57").Value = 45985.805567129632
$ws.Range("D57").Value = 45986.334270833337

$ws.Range("E22").Select()

